# Update gh-pages output data (generated at 456a3b4)
$wb = $excel.ActiveWorkbook

# ----- Sheet "展览" (Exhibitions) -----
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 2938
$ws1.Range("F3").Value = 6408
$ws1.Range("F4").Value = 2515
$ws1.Range("F6").Value = 441
$ws1.Range("F12").Value = 7336
$ws1.Range("F13").Value = 332
$ws1.Range("F14").Value = 66
$ws1.Range("F16").Value = 244
$ws1.Range("F19").Value = 8921
$ws1.Range("F20").Value = 20
$ws1.Range("F27").Value = 103
$ws1.Range("F30").Value = 43
$ws1.Range("F33").Value = 2615
$ws1.Range("F35").Value = 85
$ws1.Range("F37").Value = 1480
$ws1.Range("F38").Value = 737
$ws1.Range("F39").Value = 3852

$ws1.Range("E43").Value = "2024.06.29 09:30-06.29 16:00"
$ws1.Range("F43").Value = 19

$ws1.Range("F44").Value = 15
$ws1.Range("F46").Value = 7
$ws1.Range("F48").Value = 29
$ws1.Range("F49").Value = 42

# ----- Sheet "演出" (Performances) -----
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F5").Value = 260
$ws2.Range("G9").Value = 380

# ----- Sheet "全部类型" (All types) -----
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value = 2938
$ws4.Range("F5").Value = 260
$ws4.Range("F6").Value = 6408
$ws4.Range("F7").Value = 2515
$ws4.Range("F10").Value = 443
$ws4.Range("F18").Value = 7336
$ws4.Range("F19").Value = 332
$ws4.Range("F20").Value = 66
$ws4.Range("F22").Value = 244
$ws4.Range("F24").Value = 8921
$ws4.Range("F30").Value = 103
$ws4.Range("F31").Value = 43
$ws4.Range("F34").Value = 2615
$ws4.Range("F36").Value = 85
$ws4.Range("F38").Value = 1480
$ws4.Range("F39").Value = 737
$ws4.Range("F41").Value = 3852
$ws4.Range("F48").Value = 29
$ws4.Range("F49").Value = 42

$wb.Save()
